$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14, shifting rows 14:24 down to 15:25
$ws.Rows.Item(14).Insert()

# Fill in the new row's values
$ws.Cells.Item(14, 1).Value = 14
$ws.Cells.Item(14, 2).Value = "…"
$ws.Cells.Item(14, 3).Value = "…"

# Update selection to match the target state
$ws.Range("C14").Select()
